# Pumps table: swap the slope/offset (linear fit) formulas in columns R/S
# for rows 3-5 so that they express the input-range -> output-range mapping
# the other way around, and move the selection to S10.
#
# Old:  R = (M-L)/(P-O)          S = M-(R*P)
# New:  R = (P-O)/(M-L)          S = P-(R*M)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Formula = "=(P3-O3)/(M3-L3)"
$ws.Range("S3").Formula = "=P3-(R3*M3)"

$ws.Range("R4").Formula = "=(P4-O4)/(M4-L4)"
$ws.Range("S4").Formula = "=P4-(R4*M4)"

$ws.Range("R5").Formula = "=(P5-O5)/(M5-L5)"
$ws.Range("S5").Formula = "=P5-(R5*M5)"

# Move the active selection (as recorded in the workbook view) to S10.
$ws.Range("S10").Select()
